# Add team win/loss/tie record columns to the PHI_2014 player data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new column headers "Wins", "Losses", "Ties" ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the formatting of the existing header cell (AC1) onto the new
# header cells so they pick up the same bold/bordered/centered style
# used by the rest of row 1, instead of creating brand-new style entries.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows 2-51: every player gets the team's overall 2014 record ---
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 30).Value = 73
    $ws.Cells.Item($r, 31).Value = 89
    $ws.Cells.Item($r, 32).Value = 0
}

Write-Output "done"
